$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

# Scratch cell used to smuggle purely-numeric-looking text (e.g. cadastre
# numbers) into a cell as a real text/string value rather than a number -
# assigning such a literal straight into the destination makes the engine
# store it as a numeric value instead of a shared string.
$scratch = $ws2.Range("H1")

# ---------------------------------------------------------------------
# 1. Move the 7 "New" sheet data rows (rows 2-8) onto the end of the
#    "Previously added" sheet (rows 461-467), preserving their exact
#    shared-string-backed values and the existing row-style pattern.
# ---------------------------------------------------------------------

$moveUrls = @(
    "https://www.ss.com/msg/lv/real-estate/wood/bauska-and-reg/vecsaules-pag/kffcn.html",
    "https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/cesis/cxhdf.html",
    "https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/upmalas-pag/bhhlnf.html",
    "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/dricanu-pag/jhdxd.html",
    "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/veremu-pag/bgglm.html",
    "https://www.ss.com/msg/lv/real-estate/wood/valka-and-reg/smiltene/piixf.html",
    "https://www.ss.com/msg/lv/real-estate/wood/other/hjkix.html"
)

for ($i = 0; $i -lt 7; $i++) {
    $srcRow = 2 + $i
    $dstRow = 461 + $i

    # Register the hyperlink relationship first (this also stamps its own
    # "hyperlink" style onto the cell, which gets overwritten by the
    # format-copy step right after).
    $ws1.Hyperlinks.Add($ws1.Range("A$dstRow"), $moveUrls[$i])

    # Copy formatting from the last existing data row so the new row uses
    # the same cell styles (s=3/4/2) as every other data row.
    $ws1.Range("A460:F460").Copy()
    $ws1.Range("A${dstRow}:F${dstRow}").PasteSpecial(-4122)

    # Copy the actual values (and their shared-string / numeric types)
    # straight from the source row on the "New" sheet.
    $ws2.Range("A${srcRow}:F${srcRow}").Copy()
    $ws1.Range("A${dstRow}:F${dstRow}").PasteSpecial(-4163)
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Clear out the "New" sheet's old hyperlinks + data rows.
# ---------------------------------------------------------------------

$ws2.Cells.Hyperlinks.Delete()
$ws2.Rows("2:8").Delete()

# ---------------------------------------------------------------------
# 3. Populate the "New" sheet with the 3 freshly scraped listings.
# ---------------------------------------------------------------------

$newRows = @(
    @{ A = "https://www.ss.com/msg/lv/real-estate/wood/aizkraukle-and-reg/daudzeses-pag/jfbgx.html";
       B = "107 000 €"; C = "Aizkraukle un raj."; D = "9.44 ha.";  E = "32500090004"; F = 46056.825 },
    @{ A = "https://www.ss.com/msg/lv/real-estate/wood/limbadzi-and-reg/aloja/copmo.html";
       B = "50 000 €";  C = "Limbaži un raj.";    D = "12.40 ha."; E = "66270040050"; F = 46057.47916666667 },
    @{ A = "https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/pelecu-pag/bxolie.html";
       B = "10 000 €";  C = "Preiļi un raj.";      D = "2.82 ha.";  E = "76560050295"; F = 46057.58125 }
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = 2 + $i
    $data = $newRows[$i]

    $ws2.Hyperlinks.Add($ws2.Range("A$row"), $data.A)

    $ws1.Range("A460:F460").Copy()
    $ws2.Range("A${row}:F${row}").PasteSpecial(-4122)

    $ws2.Range("A$row").Value = $data.A
    $ws2.Range("B$row").Value = $data.B
    $ws2.Range("C$row").Value = $data.C
    $ws2.Range("D$row").Value = $data.D

    $scratch.Value = "'" + $data.E
    $scratch.Copy()
    $ws2.Range("E$row").PasteSpecial(-4163)

    $ws2.Range("F$row").Value2 = $data.F
}

$scratch.Clear()
$excel.CutCopyMode = 0
